$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price/volume data (and the Polkadot/WrappedBTC row swap)
$changes = @(
    @{ Cell = 'D2'; Value = '62.160.91' },
    @{ Cell = 'E2'; Value = '  +0.04%  ' },
    @{ Cell = 'D3'; Value = '3.012.67' },
    @{ Cell = 'E3'; Value = '  +0.54%  ' },
    @{ Cell = 'D4'; Value = '0.999' },
    @{ Cell = 'E4'; Value = '  -0.04%  ' },
    @{ Cell = 'D5'; Value = '600.73' },
    @{ Cell = 'E5'; Value = '  +3.35%  ' },
    @{ Cell = 'D6'; Value = '147.26' },
    @{ Cell = 'E6'; Value = '  +1.28%  ' },
    @{ Cell = 'E7'; Value = '  +0.03%  ' },
    @{ Cell = 'D8'; Value = '3.008.00' },
    @{ Cell = 'E8'; Value = '  +0.34%  ' },
    @{ Cell = 'D9'; Value = '0.516' },
    @{ Cell = 'E9'; Value = '  -1.56%  ' },
    @{ Cell = 'D10'; Value = '0.150' },
    @{ Cell = 'E10'; Value = '  +1.47%  ' },
    @{ Cell = 'D11'; Value = '6.19' },
    @{ Cell = 'E11'; Value = '  +7.06%  ' },
    @{ Cell = 'D12'; Value = '0.456' },
    @{ Cell = 'E12'; Value = '  +0.69%  ' },
    @{ Cell = 'D13'; Value = '0.0000230' },
    @{ Cell = 'E13'; Value = '  +1.09%  ' },
    @{ Cell = 'D14'; Value = '34.46' },
    @{ Cell = 'E14'; Value = '  +0.43%  ' },
    @{ Cell = 'E15'; Value = '  +3.26%  ' },
    @{ Cell = 'D16'; Value = '3.507.48' },
    @{ Cell = 'E16'; Value = '  +0.60%  ' },
    @{ Cell = 'B17'; Value = 'WrappedBTC' },
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc' },
    @{ Cell = 'D17'; Value = '62.014.58' },
    @{ Cell = 'E17'; Value = '  +0.02%  ' },
    @{ Cell = 'B18'; Value = 'Polkadot' },
    @{ Cell = 'C18'; Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot' },
    @{ Cell = 'D18'; Value = '6.96' },
    @{ Cell = 'E18'; Value = '  -1.51%  ' },
    @{ Cell = 'D19'; Value = '3.008.03' },
    @{ Cell = 'E19'; Value = '  +0.54%  ' },
    @{ Cell = 'D20'; Value = '449.43' },
    @{ Cell = 'E20'; Value = '  -2.21%  ' },
    @{ Cell = 'D21'; Value = '14.09' },
    @{ Cell = 'E21'; Value = '  +1.21%  ' },
    @{ Cell = 'D22'; Value = '0.687' },
    @{ Cell = 'E22'; Value = '  +0.55%  ' },
    @{ Cell = 'D23'; Value = '7.38' },
    @{ Cell = 'E23'; Value = '  -0.72%  ' },
    @{ Cell = 'D24'; Value = '81.60' },
    @{ Cell = 'E24'; Value = '  +0.50%  ' },
    @{ Cell = 'D25'; Value = '11.15' },
    @{ Cell = 'E25'; Value = '  +11.63%  ' },
    @{ Cell = 'D26'; Value = '2.24' },
    @{ Cell = 'E26'; Value = '  +1.30%  ' },
    @{ Cell = 'D27'; Value = '12.17' },
    @{ Cell = 'E27'; Value = '  -0.55%  ' },
    @{ Cell = 'E28'; Value = '  +0.24%  ' },
    @{ Cell = 'D29'; Value = '2.73' },
    @{ Cell = 'E29'; Value = '  +4.47%  ' },
    @{ Cell = 'D30'; Value = '0.997' },
    @{ Cell = 'E30'; Value = '  -0.14%  ' },
    @{ Cell = 'D31'; Value = '7.22' },
    @{ Cell = 'E31'; Value = '  +3.15%  ' },
    @{ Cell = 'D32'; Value = '2.08' },
    @{ Cell = 'E32'; Value = '  +0.10%  ' },
    @{ Cell = 'D33'; Value = '27.34' },
    @{ Cell = 'E33'; Value = '  -2.75%  ' },
    @{ Cell = 'D34'; Value = '0.111' },
    @{ Cell = 'E34'; Value = '  +2.87%  ' },
    @{ Cell = 'D35'; Value = '0.0₃0840' },
    @{ Cell = 'E35'; Value = '  +5.93%  ' },
    @{ Cell = 'D36'; Value = '1.03' },
    @{ Cell = 'E36'; Value = '  +0.21%  ' },
    @{ Cell = 'D37'; Value = '5.81' },
    @{ Cell = 'E37'; Value = '  +1.35%  ' },
    @{ Cell = 'D38'; Value = '50.54' },
    @{ Cell = 'E38'; Value = '  +0.55%  ' },
    @{ Cell = 'D39'; Value = '2.06' },
    @{ Cell = 'E39'; Value = '  -1.86%  ' },
    @{ Cell = 'D40'; Value = '8.98' },
    @{ Cell = 'E40'; Value = '  -1.90%  ' },
    @{ Cell = 'D41'; Value = '0.124' },
    @{ Cell = 'E41'; Value = '  +8.44%  ' },
    @{ Cell = 'D42'; Value = '2.92' },
    @{ Cell = 'E42'; Value = '  +1.83%  ' },
    @{ Cell = 'D43'; Value = '401.56' },
    @{ Cell = 'E43'; Value = '  +2.17%  ' },
    @{ Cell = 'D44'; Value = '40.78' },
    @{ Cell = 'E44'; Value = '  +11.26%  ' },
    @{ Cell = 'D45'; Value = '0.274' },
    @{ Cell = 'E45'; Value = '  +0.62%  ' },
    @{ Cell = 'D46'; Value = '0.0353' },
    @{ Cell = 'E46'; Value = '  -0.62%  ' },
    @{ Cell = 'D47'; Value = '2.717.23' },
    @{ Cell = 'E47'; Value = '  -0.11%  ' },
    @{ Cell = 'D48'; Value = '131.66' },
    @{ Cell = 'E48'; Value = '  +2.91%  ' },
    @{ Cell = 'E49'; Value = '  +0.12%  ' },
    @{ Cell = 'D50'; Value = '2.19' },
    @{ Cell = 'E50'; Value = '  +0.43%  ' },
    @{ Cell = 'E51'; Value = '  -1.14%  ' }
)

foreach ($item in $changes) {
    $c = $ws.Range($item.Cell)
    $origStyle = $c.Style
    # Force text number format so numeric-looking strings (e.g. "6.19", "0.516")
    # are not silently coerced into floating point numbers by Excel.
    $c.NumberFormat = "@"
    $c.Value = $item.Value
    $c.Style = $origStyle
}
